$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: column D changes from "empresa" to "servicio",
# column F changes from "cuenta" to "empresa"
$ws.Range("D1").Value = "servicio"
$ws.Range("F1").Value = "empresa"

# New test case data (row 7, ATC06_descargarCartolaLuz) gets
# a servicio value ("Luz") and an empresa value ("CGE")
$ws.Range("D7").Value = "Luz"
$ws.Range("F7").Value = "CGE"

# Update the active selection to reflect where the user left off editing
$ws.Range("D8").Select()
